$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row 4: was (2, AFTERNOON) -> becomes (3, MORNING)
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "MORNING"

# Append the new rows following the (n, MORNING) / (n, MIDDAY) pattern
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "MIDDAY"

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "MORNING"

$ws.Cells.Item(7, 1).Value = 4
$ws.Cells.Item(7, 2).Value = "MIDDAY"

$ws.Cells.Item(8, 1).Value = 5
$ws.Cells.Item(8, 2).Value = "MORNING"

$ws.Cells.Item(9, 1).Value = 5
$ws.Cells.Item(9, 2).Value = "MIDDAY"

# Match the numeric formatting already used by A2:A4 (integer number format)
$ws.Range("A5:A9").NumberFormat = "0"

# Update selection to reflect the new active cell (A9) as in the edited workbook
$ws.Range("A9").Select()
